$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: free exploration note 1
$ws.Range("A17").Value = "-"
$ws.Range("C17").Value = "In “Berichte Studiengebühren” sind keine Berichte?"
$ws.Range("D17").Value = 1

# Row 18: free exploration note 2
$ws.Range("A18").Value = "-"
$ws.Range("C18").Value = "“Newsletter” in “Newsletter (Rzettel)” umbenennen wegen Branding"
$ws.Range("D18").Value = 2

# Move the active selection to C19 (next empty row), matching author's edit
[void]$ws.Range("C19").Select()
